# v3.0 update FCI 27/1/2023
# Reorders rows so fund rows come first (sorted), followed by avg/total,
# and adds a new column C with the 13-01-2023 figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header date in column C, same style as B1
$ws.Range("B1").Copy($ws.Range("C1"))
$ws.Range("C1").Value = "13-01-2023"

# Row labels (A) and values (B = 06-01-2023, C = 13-01-2023) after reorder
$ws.Range("A2").Value = "Consultatio Renta Variable"
$ws.Range("B2").Value = 45696.49
$ws.Range("C2").Value = 46285.36

$ws.Range("A3").Value = "Delta Recursos Naturales"
$ws.Range("B3").Value = 340021.01
$ws.Range("C3").Value = 339657.57

$ws.Range("A4").Value = "Delta Select"
$ws.Range("B4").Value = 47985.63
$ws.Range("C4").Value = 47262.34

$ws.Range("A5").Value = "avg"
$ws.Range("B5").Value = 144567.71
$ws.Range("C5").Value = 144401.76

$ws.Range("A6").Value = "total"
$ws.Range("B6").Value = 433703.13
$ws.Range("C6").Value = 433205.27
